$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "License Information" heading paragraph entirely (it is
#    merged away; the body paragraph right after it keeps its own plain
#    <w:bidi/> paragraph properties, no heading style).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("License Information")
if ($found) {
    [void]$rng.Expand(4)
    $rng.Delete()
}

# ---------------------------------------------------------------------------
# 2) Replace the license-body paragraph's text (the one that used to start
#    with the bold "<Arabic Key Terms title>" run, followed by the
#    "(Arabic) is based on..." sentence, hyperlinks, etc.) with the new
#    resource-description text. Only the leading
#    "Biblica Study Notes (Key Terms)" segment stays bold.
# ---------------------------------------------------------------------------
$bodyFind = $d.Content
$foundBody = $bodyFind.Find.Execute("Biblica Bible Dictionary")
if ($foundBody) {
    [void]$bodyFind.Expand(4)
    $pStart = $bodyFind.Start
    $pEnd = $bodyFind.End - 1

    $boldText = "Biblica Study Notes (Key Terms)"
    $restText = " " + [char]0x00A9 + " 2023 Biblica Inc. Released under CC BY-SA 4.0 license. " + `
        "Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (" + `
        [char]0x0639 + [char]0x0631 + [char]0x0628 + [char]0x064A + `
        "), French (Fran" + [char]0x00E7 + "ais), Hindi (" + `
        [char]0x0939 + [char]0x093F + [char]0x0902 + [char]0x0926 + [char]0x0940 + `
        "), Indonesian (Bahasa Indonesia), Portuguese (Portugu" + [char]0x00EA + "s), Russian (" + `
        [char]0x0420 + [char]0x0443 + [char]0x0441 + [char]0x0441 + [char]0x043A + [char]0x0438 + [char]0x0439 + `
        "), Spanish (Espa" + [char]0x00F1 + "ol), Swahili (Kiswahili), and Simplified Chinese (" + `
        [char]0x7B80 + [char]0x4F53 + [char]0x4E2D + [char]0x6587 + `
        ")from Biblica Study Notes " + [char]0x00A9 + " 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."

    $fullText = $boldText + $restText

    $bodyRange = $d.Range($pStart, $pEnd)
    $bodyRange.Text = $fullText

    $boldRange = $d.Range($pStart, $pStart + $boldText.Length)
    $boldRange.Font.Bold = 1

    $restRange = $d.Range($pStart + $boldText.Length, $pStart + $fullText.Length)
    $restRange.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 3) Remove the "This PDF version is provided under the same license."
#    paragraph entirely.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("This PDF version is provided under the same license.")
if ($found3) {
    [void]$rng3.Expand(4)
    $rng3.Delete()
}

# ---------------------------------------------------------------------------
# 4) Remove the paragraph holding the italic list of key terms that used to
#    follow the single-letter "<heading>" (Heading2, 2-character) paragraph
#    in the glossary index section (right before the " " spacer paragraph).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $styleName = $p.Style.NameLocal
    if ($styleName -eq "Heading 2" -and $p.Range.Text.Length -eq 2) {
        $nextPara = $p.Next()
        if ($nextPara -ne $null) {
            $nextPara.Range.Delete()
        }
        break
    }
}
